$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerate the K column (col G) values for rows 2-5
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 2
